$d = $word.ActiveDocument

function Clear-CellText($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $full = $r.Text
    if ($full.Length -gt 2) {
        $content = $full.Substring(0, $full.Length - 2)
        if ($content.Length -gt 0) {
            $r.Find.Execute($content, $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
        }
    }
}

# ISSUE_TABLE (Table 5)
$t5 = $d.Tables.Item(5)
Clear-CellText $t5 2 1
Clear-CellText $t5 3 1
Clear-CellText $t5 4 1
Clear-CellText $t5 8 1
Clear-CellText $t5 9 1
Clear-CellText $t5 10 1

# ISSUE_PENDING_TABLE (Table 7)
$t7 = $d.Tables.Item(7)
Clear-CellText $t7 2 2
Clear-CellText $t7 2 5
Clear-CellText $t7 3 2
Clear-CellText $t7 3 5
Clear-CellText $t7 4 2
Clear-CellText $t7 4 5
Clear-CellText $t7 5 5
Clear-CellText $t7 8 2
Clear-CellText $t7 8 5
Clear-CellText $t7 9 2
Clear-CellText $t7 9 5
Clear-CellText $t7 10 2
Clear-CellText $t7 10 5
Clear-CellText $t7 11 5

# ISSUE_COVERED_TABLE (Table 9)
$t9 = $d.Tables.Item(9)
Clear-CellText $t9 2 2
Clear-CellText $t9 3 2
Clear-CellText $t9 4 2
Clear-CellText $t9 5 2
Clear-CellText $t9 8 2
Clear-CellText $t9 9 2
Clear-CellText $t9 10 2
Clear-CellText $t9 11 2
